$wb = $excel.ActiveWorkbook

# Overview sheet: mark c91ca27e-...md as handed back (row 3)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# zh-cn sheet: update status, handback datetime, and clear error detail for row 3
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-08-20 04:51:25"
$wsZhCn.Range("P3").Value = ""

# de-de sheet: update status, handback datetime, and clear error detail for row 3
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-08-20 04:51:31"
$wsDeDe.Range("P3").Value = ""
